$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for cells whose content looks numeric (e.g. "1.00", "14.00")
# by switching the data range to Text format before writing, then restoring the
# default ("Normal") style once all values are in place.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '89.102.27'
$ws.Range("E2").Value = '  +2.37%  '
$ws.Range("D3").Value = '3.273.59'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '211.76'
$ws.Range("E5").Value = '  -3.24%  '
$ws.Range("D6").Value = '625.66'
$ws.Range("E6").Value = '  -1.85%  '
$ws.Range("D7").Value = '0.374'
$ws.Range("E7").Value = '  +17.00%  '
$ws.Range("D8").Value = '0.710'
$ws.Range("E8").Value = '  +15.33%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Value = '3.273.41'
$ws.Range("E10").Value = '  -1.93%  '
$ws.Range("D11").Value = '0.571'
$ws.Range("E11").Value = '  -6.09%  '
$ws.Range("D12").Value = '0.187'
$ws.Range("E12").Value = '  +12.03%  '
$ws.Range("D13").Value = '0.0000260'
$ws.Range("E13").Value = '  -5.08%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.878.97'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").Value = '5.43'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '33.92'
$ws.Range("E16").Value = '  -0.99%  '
$ws.Range("D17").Value = '88.933.09'
$ws.Range("E17").Value = '  +2.53%  '
$ws.Range("D18").Value = '3.280.63'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("B19").Value = 'SuiNetwork'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").Value = '3.10'
$ws.Range("E19").Value = '  -4.29%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '14.00'
$ws.Range("E20").Value = '  -4.55%  '
$ws.Range("D21").Value = '435.78'
$ws.Range("E21").Value = '  -2.81%  '
$ws.Range("D22").Value = '8.83'
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("D23").Value = '5.31'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = '7.42'
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").Value = '5.21'
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("D26").Value = '12.16'
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("D27").Value = '3.461.92'
$ws.Range("E27").Value = '  -1.32%  '
$ws.Range("B28").Value = 'Litecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D28").Value = '76.66'
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0000135'
$ws.Range("E29").Value = '  +3.52%  '
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").Value = '0.180'
$ws.Range("E31").Value = '  -4.00%  '
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").Value = '8.80'
$ws.Range("E33").Value = '  -5.58%  '
$ws.Range("D34").Value = '560.11'
$ws.Range("E34").Value = '  -7.21%  '
$ws.Range("D35").Value = '1.36'
$ws.Range("E35").Value = '  -12.12%  '
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").Value = '1.95'
$ws.Range("E36").Value = '  -5.01%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").Value = '7.05'
$ws.Range("E37").Value = '  +6.79%  '
$ws.Range("D38").Value = '0.140'
$ws.Range("E38").Value = '  -7.73%  '
$ws.Range("D39").Value = '22.65'
$ws.Range("E39").Value = '  -3.54%  '
$ws.Range("D40").Value = '21.84'
$ws.Range("E40").Value = '  +2.25%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '3.08'
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("D43").Value = '0.400'
$ws.Range("E43").Value = '  -4.67%  '
$ws.Range("D44").Value = '2.02'
$ws.Range("E44").Value = '  -1.86%  '
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").Value = '155.64'
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").Value = '180.10'
$ws.Range("E47").Value = '  -5.33%  '
$ws.Range("D48").Value = '44.84'
$ws.Range("E48").Value = '  -1.91%  '
$ws.Range("D49").Value = '0.133'
$ws.Range("E49").Value = '  +17.79%  '
$ws.Range("D50").Value = '1.31'
$ws.Range("E50").Value = '  -4.96%  '
$ws.Range("D51").Value = '4.21'
$ws.Range("E51").Value = '  -1.99%  '

$dataRange.Style = "Normal"
